# Update employee absence data rows 2-11 on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 31488
$ws.Cells.Item(2, 2).Value = "Raquel Leão"
$ws.Cells.Item(2, 3).Value = "Engenharia"
$ws.Cells.Item(2, 5).Value = 7
$ws.Cells.Item(2, 6).Value = 45101
$ws.Cells.Item(2, 7).Value = 4369.6

# Row 3
$ws.Cells.Item(3, 1).Value = 88896
$ws.Cells.Item(3, 2).Value = "Lorena Moreira"
$ws.Cells.Item(3, 3).Value = "Vendas"
$ws.Cells.Item(3, 4).Value = "Doenca"
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 45106
$ws.Cells.Item(3, 7).Value = 2090.28

# Row 4
$ws.Cells.Item(4, 1).Value = 92504
$ws.Cells.Item(4, 2).Value = "Lívia Vieira"
$ws.Cells.Item(4, 3).Value = "Recursos Humanos"
$ws.Cells.Item(4, 4).Value = "Consulta medica"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 45098
$ws.Cells.Item(4, 7).Value = 4482.88

# Row 5
$ws.Cells.Item(5, 1).Value = 67746
$ws.Cells.Item(5, 2).Value = "Maria Júlia Ribeiro"
$ws.Cells.Item(5, 3).Value = "Marketing"
$ws.Cells.Item(5, 6).Value = 45102
$ws.Cells.Item(5, 7).Value = 9003.51

# Row 6
$ws.Cells.Item(6, 1).Value = 68212
$ws.Cells.Item(6, 2).Value = "Maysa Carvalho"
$ws.Cells.Item(6, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(6, 4).Value = "Consulta medica"
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 45101
$ws.Cells.Item(6, 7).Value = 5207.13

# Row 7
$ws.Cells.Item(7, 1).Value = 11744
$ws.Cells.Item(7, 2).Value = "Dr. João Vitor Jesus"
$ws.Cells.Item(7, 3).Value = "Marketing"
$ws.Cells.Item(7, 4).Value = "Consulta medica"
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 45102
$ws.Cells.Item(7, 7).Value = 9752.83

# Row 8
$ws.Cells.Item(8, 1).Value = 54530
$ws.Cells.Item(8, 2).Value = "Jade Dias"
$ws.Cells.Item(8, 3).Value = "Recursos Humanos"
$ws.Cells.Item(8, 4).Value = "Outros"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 45088
$ws.Cells.Item(8, 7).Value = 3585.17

# Row 9
$ws.Cells.Item(9, 1).Value = 39042
$ws.Cells.Item(9, 2).Value = "Julia Dias"
$ws.Cells.Item(9, 3).Value = "Engenharia"
$ws.Cells.Item(9, 4).Value = "Outros"
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 45082
$ws.Cells.Item(9, 7).Value = 7237.12

# Row 10
$ws.Cells.Item(10, 1).Value = 17008
$ws.Cells.Item(10, 2).Value = "Marcelo Correia"
$ws.Cells.Item(10, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 45079
$ws.Cells.Item(10, 7).Value = 5292.73

# Row 11
$ws.Cells.Item(11, 1).Value = 45226
$ws.Cells.Item(11, 2).Value = "Rhavi Vasconcelos"
$ws.Cells.Item(11, 3).Value = "Vendas"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45090
$ws.Cells.Item(11, 7).Value = 6073.2
